$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold font, border, centered alignment) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I and J, rows 2-38
$data = @{
    2  = @(7, 8)
    3  = @(7, 7)
    4  = @(7, 7)
    5  = @(9, 9)
    6  = @(7, 7)
    7  = @(7, 7)
    8  = @(9, 9)
    9  = @(7, 7)
    10 = @(7, 8)
    11 = @(7, 7)
    12 = @(7, 7)
    13 = @(7, 7)
    14 = @(8, 8)
    15 = @(7, 7)
    16 = @(7, 7)
    17 = @(9, 9)
    18 = @(8, 9)
    19 = @(7, 7)
    20 = @(8, 9)
    21 = @(8, 8)
    22 = @(7, 7)
    23 = @(6, 8)
    24 = @(8, 8)
    25 = @(9, 9)
    26 = @(10, 10)
    27 = @(9, 9)
    28 = @(5, 5)
    29 = @(7, 7)
    30 = @(9, 9)
    31 = @(8, 9)
    32 = @(5, 6)
    33 = @(8, 8)
    34 = @(5, 5)
    35 = @(3, 3)
    36 = @(8, 8)
    37 = @(4, 4)
    38 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
